$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 9,20
$data[0,0] = "ECs"
$data[0,1] = "Pdgfc"
$data[0,2] = "Pdgfra"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.2663606666666666
$data[0,7] = 0.799082
$data[0,8] = 0.0257989900554292
$data[0,9] = 0.0257989900554292
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.2781686666666667
$data[0,13] = 0.834506
$data[0,14] = 0.001228014730390642
$data[0,15] = 0.001228014730390642
$data[0,16] = 0.0740931914991111
$data[0,17] = 0.6668387234919999
$data[0,18] = 0.00003168153981726876
$data[0,19] = 0.00003168153981726876
$data[1,0] = "ECs"
$data[1,1] = "Pdgfc"
$data[1,2] = "Pdgfra"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.2663606666666666
$data[1,7] = 0.799082
$data[1,8] = 0.0257989900554292
$data[1,9] = 0.0257989900554292
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 225.778076
$data[1,13] = 677.3342279999999
$data[1,14] = 0.9967290940769435
$data[1,15] = 0.9967290940769435
$data[1,16] = 60.13839884207732
$data[1,17] = 541.2455895786959
$data[1,18] = 0.02571460398604803
$data[1,19] = 0.02571460398604803
$data[2,0] = "ECs"
$data[2,1] = "Pdgfc"
$data[2,2] = "Pdgfra"
$data[2,3] = "sCs"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.2663606666666666
$data[2,7] = 0.799082
$data[2,8] = 0.0257989900554292
$data[2,9] = 0.0257989900554292
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.4627536666666667
$data[2,13] = 1.388261
$data[2,14] = 0.002042891192665893
$data[2,15] = 0.002042891192665893
$data[2,16] = 0.1232593751557778
$data[2,17] = 1.109334376402
$data[2,18] = 0.00005270452956391128
$data[2,19] = 0.00005270452956391127
$data[3,0] = "FAPs"
$data[3,1] = "Pdgfc"
$data[3,2] = "Pdgfra"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.441874
$data[3,7] = 7.325622
$data[3,8] = 0.2365134606058369
$data[3,9] = 0.236513460605837
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.2781686666666667
$data[3,13] = 0.834506
$data[3,14] = 0.001228014730390642
$data[3,15] = 0.001228014730390642
$data[3,16] = 0.679252834748
$data[3,17] = 6.113275512732
$data[3,18] = 0.0002904420135596347
$data[3,19] = 0.0002904420135596347
$data[4,0] = "FAPs"
$data[4,1] = "Pdgfc"
$data[4,2] = "Pdgfra"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.441874
$data[4,7] = 7.325622
$data[4,8] = 0.2365134606058369
$data[4,9] = 0.236513460605837
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 225.778076
$data[4,13] = 677.3342279999999
$data[4,14] = 0.9967290940769435
$data[4,15] = 0.9967290940769435
$data[4,16] = 551.3216135544239
$data[4,17] = 4961.894521989816
$data[4,18] = 0.2357398473266587
$data[4,19] = 0.2357398473266588
$data[5,0] = "FAPs"
$data[5,1] = "Pdgfc"
$data[5,2] = "Pdgfra"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.441874
$data[5,7] = 7.325622
$data[5,8] = 0.2365134606058369
$data[5,9] = 0.236513460605837
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.4627536666666667
$data[5,13] = 1.388261
$data[5,14] = 0.002042891192665893
$data[5,15] = 0.002042891192665893
$data[5,16] = 1.129986147038
$data[5,17] = 10.169875323342
$data[5,18] = 0.0004831712656185959
$data[5,19] = 0.0004831712656185959
$data[6,0] = "sCs"
$data[6,1] = "Pdgfc"
$data[6,2] = "Pdgfra"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 7.616226333333333
$data[6,7] = 22.848679
$data[6,8] = 0.7376875493387338
$data[6,9] = 0.7376875493387338
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.2781686666666667
$data[6,13] = 0.834506
$data[6,14] = 0.001228014730390642
$data[6,15] = 0.001228014730390642
$data[6,16] = 2.118595524174889
$data[6,17] = 19.067359717574
$data[6,18] = 0.000905891177013739
$data[6,19] = 0.000905891177013739
$data[7,0] = "sCs"
$data[7,1] = "Pdgfc"
$data[7,2] = "Pdgfra"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 7.616226333333333
$data[7,7] = 22.848679
$data[7,8] = 0.7376875493387338
$data[7,9] = 0.7376875493387338
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 225.778076
$data[7,13] = 677.3342279999999
$data[7,14] = 0.9967290940769435
$data[7,15] = 0.9967290940769435
$data[7,16] = 1719.576927920534
$data[7,17] = 15476.19235128481
$data[7,18] = 0.7352746427642367
$data[7,19] = 0.7352746427642367
$data[8,0] = "sCs"
$data[8,1] = "Pdgfc"
$data[8,2] = "Pdgfra"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 7.616226333333333
$data[8,7] = 22.848679
$data[8,8] = 0.7376875493387338
$data[8,9] = 0.7376875493387338
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.4627536666666667
$data[8,13] = 1.388261
$data[8,14] = 0.002042891192665893
$data[8,15] = 0.002042891192665893
$data[8,16] = 3.524436661913222
$data[8,17] = 31.71992995721899
$data[8,18] = 0.001507015397483386
$data[8,19] = 0.001507015397483385

$ws.Range("A2:T10").Value = $data

